$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Inputs block
$ws.Range("A1").Value = "Inputs"
$ws.Range("B1").Value = "Values"

$ws.Range("A2").Value = "Desired P(block)"
Set-TextValue "B2" "0.01"

$ws.Range("A3").Value = "Arrival Rate"
Set-TextValue "B3" "10.0"

$ws.Range("A4").Value = "Service Rate"
Set-TextValue "B4" "12.0"

# Results block
$ws.Range("A6").Value = "Results"

$ws.Range("A7").Value = "Number of Servers"
Set-TextValue "B7" "4"

$ws.Range("A8").Value = "Actual P(block)"
Set-TextValue "B8" "0.008747498215510364"
